# "completed quotation creation functional flow"
#
# On the customer-quotation entry (Sheet1, row 2):
#   - the NIC number had a typo: 865360920V -> corrected to 865361920V
#   - the street field (column L) had been left blank and is now filled
#     in with "Rajapihilla"
#   - the on-screen view is moved on to reflect the next field of interest
#     (scrolled one column right, with the NIC cell now selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# street (column L) - previously empty, now populated
$ws.Range("L2").Value = "Rajapihilla"

# nic (column F) - correct the typo'd value
$ws.Range("F2").Value = "865361920V"

# Reflect the updated view state: window scrolled so column E is
# left-most, with F2 (the corrected nic) active/selected.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("F2").Select()
